# Sync attendance_reports: swap the order of names in the "Recorded By"
# column (G) from "dnasr281@gmail.com, System" to
# "System, dnasr281@gmail.com" wherever that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
